# Update "Elapsed Duration(Hrs)" values (column G) across the R1..R6 sheets.
# Each outage's elapsed duration string is bumped forward (the workbook was
# re-generated later, so the HH:MM:SS-style duration text increases).

$wb = $excel.ActiveWorkbook

# Sheet "R1"
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3929:43:47"
$ws.Range("G3").Value = "69:16:25"

# Sheet "R2"
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12111:07:28"
$ws.Range("G3").Value = "3240:50:57"
$ws.Range("G4").Value = "479:02:31"

# Sheet "R4"
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2956:57:17"
$ws.Range("G3").Value = "184:09:32"
$ws.Range("G4").Value = "72:21:57"
$ws.Range("G5").Value = "69:59:30"

# Sheet "R5"
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "430:56:16"

# Sheet "R6"
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "71:28:34"
